$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of E3 and F3 (E3: 0 -> 4800, F3: 4800 -> 0)
$ws.Range("E3").Value = 4800
$ws.Range("F3").Value = 0

# Update the active cell selection from H8 to I8
$ws.Range("I8").Select()
